# Sync non-localizable rule data:
# The "BannedPaths" rule (row 35) is removed from its old position among the
# Blocker-severity Bug rules, renamed to the singular "BannedPath", its
# severity is changed from Blocker to Critical, its tag is cleared, and it is
# re-inserted further down the table (ending up at row 40) among the
# Critical-severity Bug rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "BannedPaths" row (row 35); everything below shifts up one row.
$ws.Rows(35).Delete()

# Make room for the rule's new position (now row 40 after the shift above).
$ws.Rows(40).Insert()

# Populate the relocated / renamed rule.
$ws.Cells.Item(40, 1).Value = "BannedPath"
$ws.Cells.Item(40, 2).Value = "Customer packages should not install content under /libs"
$ws.Cells.Item(40, 3).Value = "Bug"
$ws.Cells.Item(40, 4).Value = "Critical"

# Update the sheet's active selection to match the author's saved state.
$null = $ws.Range("A37").Select()
